$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the two new "S/m" rows -------------------------
# Target layout (rows 27-38):
#   27 Tsensor_C                              (unchanged)
#   28 sigma_out (conductivity S/m)           NEW
#   29 sigma_out (conductivity uS/cm)         was row 28
#   30 NIST sensor reading (conductivity S/m) NEW
#   31 NIST sensor reading (conductivity uS/cm) was row 29
#   32 Percent difference                     was row 30
#   33 pure: 1.2 uS/cm                        was row 31
#   34 DI water in cup: 1.00 uS/cm            was row 32
#   36 Calibration 1                          was row 34
#   37 Tsensor_C                              was row 35
#   38 sigma_out                              was row 36
$ws.Rows(28).Insert()
$ws.Rows(30).Insert()

# --- Step 2: relabel the rows that used to be "sigma_out" / "NIST sensor reading" ---
$ws.Range("A29").Value = "sigma_out (conductivity uS/cm)"
$ws.Range("A31").Value = "NIST sensor reading (conductivity uS/cm)"

# --- Step 3: populate new row 28: sigma_out (conductivity S/m) ------------
$ws.Range("A28").Value = "sigma_out (conductivity S/m)"
$ws.Range("B28").Value = ""
$ws.Range("C28:I28").Formula = "=C29*10^-4"
$ws.Rows(28).RowHeight = 33

# --- Step 4: populate new row 30: NIST sensor reading (conductivity S/m) --
$ws.Range("A30").Value = "NIST sensor reading (conductivity S/m)"
$ws.Range("B30").Value = ""
$ws.Range("C30:I30").Formula = "=C31*10^-4"
$ws.Rows(30).RowHeight = 41.25

# --- Step 5: merge the label cells of the new rows -------------------------
$ws.Range("A28:B28").Merge()
$ws.Range("A30:B30").Merge()

# --- Step 6: formatting for the new rows -----------------------------------
# Row 28 data cells (C:I) get the same look as the Tsensor_C row (27): plain
# font size 14, border on left/right/top only (no bottom) - matches style 6.
$ws.Range("C28:I28").Font.Bold = $false
$ws.Range("C28:I28").Font.Size = 14
$ws.Range("C28:I28").HorizontalAlignment = -4108
$ws.Range("C28:I28").VerticalAlignment = -4108
$ws.Range("C28:I28").Borders.Item(7).LineStyle = 1
$ws.Range("C28:I28").Borders.Item(8).LineStyle = 1
$ws.Range("C28:I28").Borders.Item(10).LineStyle = 1
$ws.Range("C28:I28").Borders.Item(9).LineStyle = -4142

# Row 30 data cells (C:I) look like the old sigma_out row (now 29): bold
# font size 14, border all around - matches style 3.
$ws.Range("C30:I30").Font.Bold = $true
$ws.Range("C30:I30").Font.Size = 14
$ws.Range("C30:I30").HorizontalAlignment = -4108
$ws.Range("C30:I30").VerticalAlignment = -4108
$ws.Range("C30:I30").Borders.LineStyle = 1

# Label cells (A/B) for both new rows: bold size-14 font, centered, with a
# border framing the merged label (left+top+bottom on A, right+top+bottom on B).
$labelRanges = @("A28", "B28", "A30", "B30")
foreach ($addr in $labelRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Font.Size = 14
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}
$ws.Range("A28").Borders.Item(7).LineStyle = 1
$ws.Range("A28").Borders.Item(8).LineStyle = 1
$ws.Range("A28").Borders.Item(9).LineStyle = 1
$ws.Range("A28").Borders.Item(10).LineStyle = -4142
$ws.Range("B28").Borders.Item(10).LineStyle = 1
$ws.Range("B28").Borders.Item(8).LineStyle = 1
$ws.Range("B28").Borders.Item(9).LineStyle = 1
$ws.Range("B28").Borders.Item(7).LineStyle = -4142

$ws.Range("A30").Borders.Item(7).LineStyle = 1
$ws.Range("A30").Borders.Item(8).LineStyle = 1
$ws.Range("A30").Borders.Item(9).LineStyle = 1
$ws.Range("A30").Borders.Item(10).LineStyle = -4142
$ws.Range("B30").Borders.Item(10).LineStyle = 1
$ws.Range("B30").Borders.Item(8).LineStyle = 1
$ws.Range("B30").Borders.Item(9).LineStyle = 1
$ws.Range("B30").Borders.Item(7).LineStyle = -4142

Write-Output "edit complete"
